# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2  = @{ B = 3.286832544864788;    C = 1.655778082260271;  D = 0.7527432677738641;  E = 0.4942365360607697;  G = 6.189590430959694 }
    3  = @{ B = 1.455362044514542;    C = 1.655778082260271;  D = 3.537761648806719;   E = 0.4942365360607697;  G = 7.143138311642302 }
    4  = @{ B = 3.286832544864788;    C = 1.655778082260271;  D = 0.7527432677738641;  E = 0.4942365360607697;  G = 6.189590430959694 }
    5  = @{ B = 0.01293466051926884; C = 0.04071648406533734; D = 0.1494219747398047;  E = 0.4942365360607697;  G = 0.6973096553851805 }
    6  = @{ B = 3.286832544864788;    C = 1.655778082260271;  D = 0.1494219747398047;  E = 0.4942365360607697;  G = 5.586269137925634 }
    7  = @{ B = 0.0006408296065709695; C = 10.34677158129881; D = 3.537761648806719;   E = 10.19245300693656;   G = 24.07762706664866 }
    8  = @{ B = 3.286832544864788;    C = 1.655778082260271;  D = 0.7527432677738641;  E = 0.4942365360607697;  G = 6.189590430959694 }
    9  = @{ B = 1.455362044514542;    C = 1.655778082260271;  D = 0.7527432677738641;  E = 0.4942365360607697;  G = 4.358119930609447 }
    10 = @{ B = 0.6606524410359556;   C = 1.655778082260271;  D = 0.1494219747398047;  E = 0.4942365360607697;  G = 2.960089034096801 }
    11 = @{ B = 0.1190320826869504;   C = 0.306821227259698;  D = 3.537761648806719;   E = 10.19245300693656;   G = 14.15606796568992 }
    12 = @{ B = 0.1190320826869504;   C = 0.306821227259698;  D = 0.1494219747398047;  E = 0.4942365360607697;  G = 1.069511820747223 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
